$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns: D = CanRaid, E = IsArmy
$ws.Range("D1").Value2 = "CanRaid"
$ws.Range("E1").Value2 = "IsArmy"

# CanRaid = 1 for: sand, silicon, carbon, iron, chip, deuter, floatmod
$ws.Range("D2").Value2 = 1   # sand
$ws.Range("D3").Value2 = 1   # silicon
$ws.Range("D5").Value2 = 1   # carbon
$ws.Range("D6").Value2 = 1   # iron
$ws.Range("D7").Value2 = 1   # chip
$ws.Range("D11").Value2 = 1  # deuter
$ws.Range("D15").Value2 = 1  # floatmod

# IsArmy = 1 for: tank, chopper, ship
$ws.Range("E8").Value2 = 1   # tank
$ws.Range("E9").Value2 = 1   # chopper
$ws.Range("E10").Value2 = 1  # ship

# Page setup matching the target (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Match the final selected cell/range as recorded in the saved workbook
$null = $ws.Range("E10").Select()

$wb.Save()
